$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.918.74"
$ws.Range("E2").Value = "  -0.57%  "

# Row 3
$ws.Range("D3").Value = "1.816.74"
$ws.Range("E3").Value = "  +0.21%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'309.91"
$ws.Range("E5").Value = "  -0.48%  "

# Row 6
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.16%  "

# Row 7
$ws.Range("D7").Value = "'0.4666"
$ws.Range("E7").Value = "  +0.86%  "

# Row 8
$ws.Range("D8").Value = "'0.3699"
$ws.Range("E8").Value = "  -1.59%  "

# Row 9
$ws.Range("D9").Value = "'0.07376"
$ws.Range("E9").Value = "  -0.43%  "

# Row 10
$ws.Range("D10").Value = "'0.8719"
$ws.Range("E10").Value = "  +0.97%  "

# Row 11
$ws.Range("D11").Value = "'20.45"

# Row 12
$ws.Range("D12").Value = "1.866.28"
$ws.Range("E12").Value = "  +2.95%  "

# Row 13
$ws.Range("D13").Value = "'5.377"
$ws.Range("E13").Value = "  -0.25%  "

# Row 14
$ws.Range("D14").Value = "'92.45"
$ws.Range("E14").Value = "  +0.67%  "

# Row 15
$ws.Range("D15").Value = "'0.07082"
$ws.Range("E15").Value = "  +0.14%  "

# Row 16
$ws.Range("D16").Value = "'6.512"
$ws.Range("E16").Value = "  -2.05%  "

# Row 17
$ws.Range("E17").Value = "  -0.12%  "

# Row 18
$ws.Range("D18").Value = "'0.000008720"

# Row 19
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("D20").Value = "'14.76"
$ws.Range("E20").Value = "  -0.84%  "

# Row 21
$ws.Range("D21").Value = "26.939.88"
$ws.Range("E21").Value = "  -0.51%  "

# Row 22
$ws.Range("D22").Value = "'5.342"
$ws.Range("E22").Value = "  +0.23%  "

# Row 23
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "  -2.73%  "

# Row 24
$ws.Range("D24").Value = "2.040.54"
$ws.Range("E24").Value = "  -0.04%  "

# Row 25
$ws.Range("D25").Value = "'1.901"
$ws.Range("E25").Value = "  -0.58%  "

# Row 26
$ws.Range("D26").Value = "'151.73"
$ws.Range("E26").Value = "  +0.13%  "

# Row 27
$ws.Range("D27").Value = "'2.207"
$ws.Range("E27").Value = "  +0.56%  "

# Row 28
$ws.Range("D28").Value = "'18.43"
$ws.Range("E28").Value = "  -0.61%  "

# Row 29
$ws.Range("D29").Value = "'5.319"
$ws.Range("E29").Value = "  +1.04%  "

# Row 30
$ws.Range("D30").Value = "'115.65"
$ws.Range("E30").Value = "  -1.11%  "

# Row 31
$ws.Range("D31").Value = "'0.08927"
$ws.Range("E31").Value = "  -0.03%  "

# Row 32
$ws.Range("D32").Value = "'0.7689"
$ws.Range("E32").Value = "  -0.42%  "

# Row 33
$ws.Range("E33").Value = "  -0.74%  "

# Row 34
$ws.Range("D34").Value = "'4.489"
$ws.Range("E34").Value = "  -0.66%  "

# Row 35
$ws.Range("D35").Value = "'2.917"
$ws.Range("E35").Value = "  +0.59%  "

# Row 36
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "  -0.09%  "

# Row 37
$ws.Range("D37").Value = "'1.093"
$ws.Range("E37").Value = "  -3.03%  "

# Row 38
$ws.Range("E38").Value = "  +0.20%  "

# Row 39
$ws.Range("D39").Value = "'0.05281"
$ws.Range("E39").Value = "  +0.71%  "

# Row 40
$ws.Range("D40").Value = "'2.966"
$ws.Range("E40").Value = "  +1.26%  "

# Row 41
$ws.Range("D41").Value = "'0.5362"
$ws.Range("E41").Value = "  +1.39%  "

# Row 42
$ws.Range("D42").Value = "'7.260"
$ws.Range("E42").Value = "  +0.33%  "

# Row 43
$ws.Range("D43").Value = "'2.377"
$ws.Range("E43").Value = "  +1.61%  "

# Row 44
$ws.Range("D44").Value = "'0.1665"
$ws.Range("E44").Value = "  -0.67%  "

# Row 45
$ws.Range("D45").Value = "'8.460"
$ws.Range("E45").Value = "  -1.72%  "

# Row 46
$ws.Range("D46").Value = "'0.4946"
$ws.Range("E46").Value = "  -1.61%  "

# Row 47
$ws.Range("D47").Value = "'10.48"
$ws.Range("E47").Value = "  +0.83%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.676"
$ws.Range("E48").Value = "  +0.21%  "

# Row 49
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  -0.08%  "

# Row 50
$ws.Range("D50").Value = "'102.75"
$ws.Range("E50").Value = "  -1.80%  "

# Row 51
$ws.Range("D51").Value = "'0.06292"
$ws.Range("E51").Value = "  -0.68%  "
